{"js": "// Insert a new \"List Bullet\" paragraph right after the\n// \"Docente(s) Respons\u00e1vel(eis) \" heading paragraph, containing the\n// professor's name/ID line: \"1814052 - Silvio Silverio da Silva\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"Docente(s) Respons\u00e1vel(eis)\" paragraph (trim to tolerate the\n// trailing space present in the source document).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text.trim();\n  if (text === \"Docente(s) Respons\u00e1vel(eis)\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find paragraph \"Docente(s) Respons\u00e1vel(eis)\"');\n}\n\n// Insert the new paragraph right after it, and apply the ListBullet style.\nconst newPara = target.insertParagraph(\n  \"1814052 - Silvio Silverio da Silva\",\n  Word.InsertLocation.after\n);\nnewPara.style = \"List Bullet\";\n\nawait context.sync();\n", "ps1": "# Insert a new \"List Bullet\" paragraph right after the\n# \"Docente(s) Respons\u00e1vel(eis) \" heading paragraph, containing the\n# professor's name/ID line: \"1814052 - Silvio Silverio da Silva\".\n\n$d = $word.ActiveDocument\n\n$targetIndex = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    $text = $p.Range.Text.Trim()\n    if ($text -eq \"Docente(s) Respons\u00e1vel(eis)\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw 'Could not find paragraph \"Docente(s) Respons\u00e1vel(eis)\"'\n}\n\n$target = $d.Paragraphs.Item($targetIndex)\n\n# Add a new empty paragraph right after the target paragraph.\n$target.Range.InsertParagraphAfter()\n\n# The freshly inserted paragraph now sits at $targetIndex + 1.\n$newP = $d.Paragraphs.Item($targetIndex + 1)\n\n$newP.Range.Text = \"1814052 - Silvio Silverio da Silva\"\n$newP.Style = \"List Bullet\"\n"}
